$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)


$ws.Range("G2").Value = 21.267222
$ws.Range("H2").Value = 63.801666
$ws.Range("I2").Value = 0.06271644651145813
$ws.Range("J2").Value = 0.06271644651145813
$ws.Range("M2").Value = 0.07423700000000001
$ws.Range("N2").Value = 0.222711
$ws.Range("O2").Value = 0.006982421219793508
$ws.Range("P2").Value = 0.006982421219793506
$ws.Range("Q2").Value = 1.578814759614
$ws.Range("R2").Value = 14.209332836526
$ws.Range("S2").Value = 0.0004379126469516498
$ws.Range("T2").Value = 0.0004379126469516497
$ws.Range("G3").Value = 21.267222
$ws.Range("H3").Value = 63.801666
$ws.Range("I3").Value = 0.06271644651145813
$ws.Range("J3").Value = 0.06271644651145813
$ws.Range("O3").Value = 0.009342124750861835
$ws.Range("P3").Value = 0.009342124750861833
$ws.Range("Q3").Value = 2.112373914224
$ws.Range("R3").Value = 19.011365228016
$ws.Range("S3").Value = 0.0005859048672407954
$ws.Range("T3").Value = 0.0005859048672407953
$ws.Range("G4").Value = 21.267222
$ws.Range("H4").Value = 63.801666
$ws.Range("I4").Value = 0.06271644651145813
$ws.Range("J4").Value = 0.06271644651145813
$ws.Range("M4").Value = 0.08773500000000001
$ws.Range("N4").Value = 0.263205
$ws.Range("O4").Value = 0.00825198655277804
$ws.Range("P4").Value = 0.00825198655277804
$ws.Range("Q4").Value = 1.86587972217
$ws.Range("R4").Value = 16.79291749953
$ws.Range("S4").Value = 0.0005175352732505758
$ws.Range("T4").Value = 0.0005175352732505758
$ws.Range("G5").Value = 21.267222
$ws.Range("H5").Value = 63.801666
$ws.Range("I5").Value = 0.06271644651145813
$ws.Range("J5").Value = 0.06271644651145813
$ws.Range("M5").Value = 10.370688
$ws.Range("N5").Value = 31.112064
$ws.Range("O5").Value = 0.9754234674765666
$ws.Range("P5").Value = 0.9754234674765666
$ws.Range("Q5").Value = 220.555723988736
$ws.Range("R5").Value = 1985.001515898624
$ws.Range("S5").Value = 0.06117509372401511
$ws.Range("T5").Value = 0.06117509372401511
$ws.Range("I6").Value = 0.4054090708715844
$ws.Range("J6").Value = 0.4054090708715843
$ws.Range("M6").Value = 0.07423700000000001
$ws.Range("N6").Value = 0.222711
$ws.Range("O6").Value = 0.006982421219793508
$ws.Range("P6").Value = 0.006982421219793506
$ws.Range("Q6").Value = 10.20570935339134
$ws.Range("R6").Value = 91.851384180522
$ws.Range("S6").Value = 0.002830736899150521
$ws.Range("T6").Value = 0.00283073689915052
$ws.Range("I7").Value = 0.4054090708715844
$ws.Range("J7").Value = 0.4054090708715843
$ws.Range("O7").Value = 0.009342124750861835
$ws.Range("P7").Value = 0.009342124750861833
$ws.Range("S7").Value = 0.003787382115213328
$ws.Range("T7").Value = 0.003787382115213327
$ws.Range("I8").Value = 0.4054090708715844
$ws.Range("J8").Value = 0.4054090708715843
$ws.Range("M8").Value = 0.08773500000000001
$ws.Range("N8").Value = 0.263205
$ws.Range("O8").Value = 0.00825198655277804
$ws.Range("P8").Value = 0.00825198655277804
$ws.Range("Q8").Value = 12.06134286299
$ws.Range("R8").Value = 108.55208576691
$ws.Range("S8").Value = 0.003345430201206554
$ws.Range("T8").Value = 0.003345430201206553
$ws.Range("I9").Value = 0.4054090708715844
$ws.Range("J9").Value = 0.4054090708715843
$ws.Range("M9").Value = 10.370688
$ws.Range("N9").Value = 31.112064
$ws.Range("O9").Value = 0.9754234674765666
$ws.Range("P9").Value = 0.9754234674765666
$ws.Range("Q9").Value = 1425.707228507392
$ws.Range("R9").Value = 12831.36505656653
$ws.Range("S9").Value = 0.395445521656014
$ws.Range("T9").Value = 0.3954455216560139
$ws.Range("G10").Value = 121.820091
$ws.Range("H10").Value = 365.460273
$ws.Range("I10").Value = 0.3592440621169263
$ws.Range("J10").Value = 0.3592440621169263
$ws.Range("M10").Value = 0.07423700000000001
$ws.Range("N10").Value = 0.222711
$ws.Range("O10").Value = 0.006982421219793508
$ws.Range("P10").Value = 0.006982421219793506
$ws.Range("Q10").Value = 9.043558095567002
$ws.Range("R10").Value = 81.39202286010301
$ws.Range("S10").Value = 0.002508393362410044
$ws.Range("T10").Value = 0.002508393362410043
$ws.Range("G11").Value = 121.820091
$ws.Range("H11").Value = 365.460273
$ws.Range("I11").Value = 0.3592440621169263
$ws.Range("J11").Value = 0.3592440621169263
$ws.Range("O11").Value = 0.009342124750861835
$ws.Range("P11").Value = 0.009342124750861833
$ws.Range("Q11").Value = 12.099821145272
$ws.Range("R11").Value = 108.898390307448
$ws.Range("S11").Value = 0.003356102844302684
$ws.Range("T11").Value = 0.003356102844302683
$ws.Range("G12").Value = 121.820091
$ws.Range("H12").Value = 365.460273
$ws.Range("I12").Value = 0.3592440621169263
$ws.Range("J12").Value = 0.3592440621169263
$ws.Range("M12").Value = 0.08773500000000001
$ws.Range("N12").Value = 0.263205
$ws.Range("O12").Value = 0.00825198655277804
$ws.Range("P12").Value = 0.00825198655277804
$ws.Range("Q12").Value = 10.687885683885
$ws.Range("R12").Value = 96.19097115496501
$ws.Range("S12").Value = 0.002964477169754235
$ws.Range("T12").Value = 0.002964477169754235
$ws.Range("G13").Value = 121.820091
$ws.Range("H13").Value = 365.460273
$ws.Range("I13").Value = 0.3592440621169263
$ws.Range("J13").Value = 0.3592440621169263
$ws.Range("M13").Value = 10.370688
$ws.Range("N13").Value = 31.112064
$ws.Range("O13").Value = 0.9754234674765666
$ws.Range("P13").Value = 0.9754234674765666
$ws.Range("Q13").Value = 1263.358155892608
$ws.Range("R13").Value = 11370.22340303347
$ws.Range("S13").Value = 0.3504150887404593
$ws.Range("T13").Value = 0.3504150887404593
$ws.Range("G14").Value = 58.539182
$ws.Range("H14").Value = 175.617546
$ws.Range("I14").Value = 0.1726304205000311
$ws.Range("J14").Value = 0.1726304205000311
$ws.Range("M14").Value = 0.07423700000000001
$ws.Range("N14").Value = 0.222711
$ws.Range("O14").Value = 0.006982421219793508
$ws.Range("P14").Value = 0.006982421219793506
$ws.Range("Q14").Value = 4.345773254134001
$ws.Range("R14").Value = 39.11195928720601
$ws.Range("S14").Value = 0.001205378311281293
$ws.Range("T14").Value = 0.001205378311281293
$ws.Range("G15").Value = 58.539182
$ws.Range("H15").Value = 175.617546
$ws.Range("I15").Value = 0.1726304205000311
$ws.Range("J15").Value = 0.1726304205000311
$ws.Range("O15").Value = 0.009342124750861835
$ws.Range("P15").Value = 0.009342124750861833
$ws.Range("Q15").Value = 5.814423765210667
$ws.Range("R15").Value = 52.32981388689601
$ws.Range("S15").Value = 0.001612734924105027
$ws.Range("T15").Value = 0.001612734924105027
$ws.Range("G16").Value = 58.539182
$ws.Range("H16").Value = 175.617546
$ws.Range("I16").Value = 0.1726304205000311
$ws.Range("J16").Value = 0.1726304205000311
$ws.Range("M16").Value = 0.08773500000000001
$ws.Range("N16").Value = 0.263205
$ws.Range("O16").Value = 0.00825198655277804
$ws.Range("P16").Value = 0.00825198655277804
$ws.Range("Q16").Value = 5.135935132770001
$ws.Range("R16").Value = 46.22341619493
$ws.Range("S16").Value = 0.001424543908566675
$ws.Range("T16").Value = 0.001424543908566675
$ws.Range("G17").Value = 58.539182
$ws.Range("H17").Value = 175.617546
$ws.Range("I17").Value = 0.1726304205000311
$ws.Range("J17").Value = 0.1726304205000311
$ws.Range("M17").Value = 10.370688
$ws.Range("N17").Value = 31.112064
$ws.Range("O17").Value = 0.9754234674765666
$ws.Range("P17").Value = 0.9754234674765666
$ws.Range("Q17").Value = 607.091592297216
$ws.Range("R17").Value = 5463.824330674945
$ws.Range("S17").Value = 0.1683877633560781
$ws.Range("T17").Value = 0.1683877633560781
